$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B20").Value = 99
